$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: add date, set comment text
$ws.Cells.Item(26, 1).Value = 44091
$ws.Cells.Item(26, 1).NumberFormat = $ws.Cells.Item(25, 1).NumberFormat
$ws.Cells.Item(26, 5).Value = "visual studio hajonnut - n. 3 tuntia korjaukseen"

# Row 27: add date, start/end times, comment text
$ws.Cells.Item(27, 1).Value = 44091
$ws.Cells.Item(27, 1).NumberFormat = $ws.Cells.Item(25, 1).NumberFormat
$ws.Cells.Item(27, 2).Value = 0.79166666666666663
$ws.Cells.Item(27, 2).NumberFormat = $ws.Cells.Item(25, 2).NumberFormat
$ws.Cells.Item(27, 3).Value = 0.95833333333333337
$ws.Cells.Item(27, 3).NumberFormat = $ws.Cells.Item(25, 3).NumberFormat
$ws.Cells.Item(27, 5).Value = "Bugfixes/Skill Increase/feature designs"

# Update the selection to E27 like in the saved workbook
$ws.Range("E27").Select() | Out-Null
